$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "2025-07-31" / 390878 / SUPORTE PARA NOTEBOOK row (was row 2),
# shifting the rest of the rows up by one.
$ws.Rows.Item(2).Delete()

# After the shift, update the cells whose values changed in the new export.
$ws.Range("G2").Value = -1325
$ws.Range("G5").Value = -1325
$ws.Range("G6").Value = -1325
$ws.Range("G7").Value = -1325
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = -93
$ws.Range("G11").Value = -53
$ws.Range("I11").Value = 0.14
$ws.Range("G12").Value = -1325
$ws.Range("G13").Value = -153
$ws.Range("I13").Value = 0.14
